$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old 3x4 test data so stale values don't linger in cells that
# are empty in the new layout.
$ws.Range("A1:R5").ClearContents()

# ---------------------------------------------------------------------------
# Header row (row 1) - bold text, bordered cells
# ---------------------------------------------------------------------------
$headers = @(
    "Description",
    "FirstName",
    "LastName",
    "UserName",
    "Password",
    "PasswordConfirm",
    "Birthday-Month",
    "Birthday-Day",
    "Birthday-Year",
    "Gender",
    "NameMesEr",
    "UserNameMesEr",
    "PasswordMesEr",
    "PasswordConfirmMesEr",
    "Birthday-MonthMesEr",
    "Birthday-DayMesEr",
    "Birthday-YearMesEr",
    "GenderMesEr"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
}

$headerRange = $ws.Range("A1:R1")
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true

# ---------------------------------------------------------------------------
# Data rows 2-5, column A descriptions + columns B-R test data
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "All field is blank"
$ws.Range("B2").Value = "   "
$ws.Range("C2").Value = " "
$ws.Range("D2").Value = " "

$ws.Range("A3").Value = "Only firstname and lastname are blank"
$ws.Range("C3").Value = "Huy"
$ws.Range("D3").Value = "yoyo"
$ws.Range("E3").Value = " "

$ws.Range("A4").Value = "Only firstname is blank"
$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 12

$ws.Range("A5").Value = "Only lastname is blank"

# Apply the thin-border style to every used cell A2:R5 (no bold)
$dataRange = $ws.Range("A2:R5")
$dataRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Column widths (values chosen so the engine's internal rounding lands as
# close as possible to the target OOXML column widths)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 35.8333333333333
$ws.Columns.Item(2).ColumnWidth = 9.33333333333333
$ws.Columns.Item(3).ColumnWidth = 8.83333333333333
$ws.Columns.Item(4).ColumnWidth = 9.66666666666667
$ws.Columns.Item(5).ColumnWidth = 8.66666666666667
$ws.Columns.Item(6).ColumnWidth = 15.8333333333333
$ws.Columns.Item(7).ColumnWidth = 14.5
$ws.Columns.Item(8).ColumnWidth = 11.6666666666667
$ws.Columns.Item(9).ColumnWidth = 12.5
$ws.Columns.Item(10).ColumnWidth = 6.83333333333333
$ws.Columns.Item(11).ColumnWidth = 11
$ws.Columns.Item(12).ColumnWidth = 15.3333333333333
$ws.Columns.Item(13).ColumnWidth = 14.3333333333333
$ws.Columns.Item(14).ColumnWidth = 21.6666666666667
$ws.Columns.Item(15).ColumnWidth = 20.1666666666667
$ws.Columns.Item(16).ColumnWidth = 17.5
$ws.Columns.Item(17).ColumnWidth = 18.1666666666667
$ws.Columns.Item(18).ColumnWidth = 12.6666666666667

# ---------------------------------------------------------------------------
# Page setup
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$ws.Range("A12").Select()
